$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3400
$ws.Range("I69").Value = 3400
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 10200
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -9326

# Row 72
$ws.Range("H72").Value = 3400
$ws.Range("I72").Value = 3400
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 30600
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -26232

# Row 88
$ws.Range("H88").Value = 41669370
$ws.Range("I88").Value = 3001
$ws.Range("J88").Value = 66669188
$ws.Range("K88").Value = 3001
$ws.Range("L88").Value = 66669188
$ws.Range("M88").Value = -2595
$ws.Range("N88").Value = -66670000

# Row 91
$ws.Range("H91").Value = 41669370
$ws.Range("I91").Value = 3001
$ws.Range("J91").Value = 66669188
$ws.Range("K91").Value = 3001
$ws.Range("L91").Value = 66669188
$ws.Range("M91").Value = -1597
$ws.Range("N91").Value = -66671996

# Row 116
$ws.Range("H116").Value = 4747.6665
$ws.Range("J116").Value = 5297.2
$ws.Range("L116").Value = 5297.2
$ws.Range("N116").Value = -12181.2

# Row 129
$ws.Range("H129").Value = 19026.582
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 22680.695
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 68042.08499999999
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -78042.08499999999

# Row 132
$ws.Range("H132").Value = 2803143.2
$ws.Range("I132").Value = 3403417.2
$ws.Range("J132").Value = 1864.1111
$ws.Range("K132").Value = 10210251.6
$ws.Range("L132").Value = 5592.3333
$ws.Range("M132").Value = -10207721.6
$ws.Range("N132").Value = -10652.3333

# Row 137
$ws.Range("H137").Value = 1485.9767
$ws.Range("I137").Value = 799.9
$ws.Range("J137").Value = 3069.2307
$ws.Range("K137").Value = 2399.7
$ws.Range("L137").Value = 9207.6921
$ws.Range("M137").Value = 150.3000000000002
$ws.Range("N137").Value = -14307.6921

# Row 141
$ws.Range("H141").Value = 2597.4
$ws.Range("I141").Value = 1671
$ws.Range("J141").Value = 4450.2
$ws.Range("K141").Value = 5013
$ws.Range("L141").Value = 13350.6
$ws.Range("M141").Value = 167
$ws.Range("N141").Value = -23710.6


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2889.13
$ws.Range("I32").Value = 2889.13
$ws.Range("K32").Value = 2889.13
$ws.Range("M32").Value = -2602.13

# Row 63
$ws.Range("H63").Value = 1667666.6
$ws.Range("I63").Value = 2000900
$ws.Range("K63").Value = 2000900
$ws.Range("M63").Value = -2000214

# Row 66
$ws.Range("H66").Value = 1667666.6
$ws.Range("I66").Value = 2000900
$ws.Range("K66").Value = 10004500
$ws.Range("M66").Value = -10001068

# Row 122
$ws.Range("H122").Value = 1871.7142
$ws.Range("I122").Value = 1874.5
$ws.Range("J122").Value = 1862.8
$ws.Range("K122").Value = 5623.5
$ws.Range("L122").Value = 5588.4
$ws.Range("M122").Value = -3173.5
$ws.Range("N122").Value = -10488.4


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2926824.2
$ws.Range("I31").Value = 2194.1162
$ws.Range("J31").Value = 11909617
$ws.Range("K31").Value = 2194.1162
$ws.Range("L31").Value = 11909617
$ws.Range("M31").Value = -1899.1162
$ws.Range("N31").Value = -11910207

# Row 34
$ws.Range("H34").Value = 2926824.2
$ws.Range("I34").Value = 2194.1162
$ws.Range("J34").Value = 11909617
$ws.Range("K34").Value = 2194.1162
$ws.Range("L34").Value = 11909617
$ws.Range("M34").Value = -1992.1162
$ws.Range("N34").Value = -11910021

# Row 132
$ws.Range("H132").Value = 5002197
$ws.Range("I132").Value = 1274.909
$ws.Range("J132").Value = 8931493
$ws.Range("K132").Value = 3824.727
$ws.Range("L132").Value = 26794479
$ws.Range("M132").Value = -1294.727
$ws.Range("N132").Value = -26799539


$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 100000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

# Row 78
$ws.Range("H78").Value = 100000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

# Row 131
$ws.Range("H131").Value = 1916.1837
$ws.Range("I131").Value = 10322.728
$ws.Range("J131").Value = 853.2873499999999
$ws.Range("K131").Value = 30968.184
$ws.Range("L131").Value = 2559.86205
$ws.Range("M131").Value = -25928.184
$ws.Range("N131").Value = -12639.86205

# Row 132
$ws.Range("H132").Value = 1325.2444
$ws.Range("I132").Value = 764.2963
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 6878.6667
$ws.Range("L132").Value = 19500.0003
$ws.Range("M132").Value = -4348.6667
$ws.Range("N132").Value = -24560.0003

# Row 141
$ws.Range("H141").Value = 2916.3845
$ws.Range("I141").Value = 2572.6667
$ws.Range("J141").Value = 4360
$ws.Range("K141").Value = 7718.000100000001
$ws.Range("L141").Value = 13080
$ws.Range("M141").Value = -2538.000100000001
$ws.Range("N141").Value = -23440


$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 95166.55
$ws.Range("I132").Value = 146632
$ws.Range("J132").Value = 5102
$ws.Range("K132").Value = 439896
$ws.Range("L132").Value = 15306
$ws.Range("M132").Value = -437366
$ws.Range("N132").Value = -20366


$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 20835676
$ws.Range("I100").Value = 33334994
$ws.Range("J100").Value = 3482.3333
$ws.Range("K100").Value = 33334994
$ws.Range("L100").Value = 3482.3333
$ws.Range("M100").Value = -33334453
$ws.Range("N100").Value = -4564.3333

# Row 132
$ws.Range("H132").Value = 5407.4893
$ws.Range("I132").Value = 7061.5312
$ws.Range("J132").Value = 1878.8667
$ws.Range("K132").Value = 21184.5936
$ws.Range("L132").Value = 5636.6001
$ws.Range("M132").Value = -18654.5936
$ws.Range("N132").Value = -10696.6001


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 730.86664
$ws.Range("I113").Value = 796.3
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 2388.9
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -218.8999999999996
$ws.Range("N113").Value = -6140

# Row 136
$ws.Range("H136").Value = 1105.2
$ws.Range("I136").Value = 1045.1951
$ws.Range("J136").Value = 1234.6842
$ws.Range("K136").Value = 3135.5853
$ws.Range("L136").Value = 3704.0526
$ws.Range("M136").Value = -585.5852999999997
$ws.Range("N136").Value = -8804.052599999999

